$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.591.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.53%  '

$ws.Range("D3").Value = "'2.434.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'514.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.14%  '

$ws.Range("D6").Value = "'129.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.84%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("E8").Value = '  -1.90%  '

$ws.Range("D9").Value = "'2.450.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.41%  '

$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("D11").Value = "'0.0953"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.67%  '

$ws.Range("D13").Value = "'0.331"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.69%  '

$ws.Range("D14").Value = "'2.868.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.89%  '

$ws.Range("D15").Value = "'57.538.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.45%  '

$ws.Range("D16").Value = "'21.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.59%  '

$ws.Range("D17").Value = "'0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.99%  '

$ws.Range("D18").Value = "'2.442.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.80%  '

$ws.Range("E19").Value = '  -4.24%  '

$ws.Range("D20").Value = "'315.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.68%  '

$ws.Range("E21").Value = '  -2.53%  '

$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").Value = "'5.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.71%  '

$ws.Range("D24").Value = "'63.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("D25").Value = "'0.407"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.60%  '

$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("E28").Value = '  -3.41%  '

$ws.Range("D29").Value = "'170.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.98%  '

$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").Value = "'6.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.57%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = "'0.0₃0723"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.76%  '

$ws.Range("E32").Value = '  -2.42%  '

$ws.Range("D33").Value = "'1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.96%  '

$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("D36").Value = "'17.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.93%  '

$ws.Range("D38").Value = "'3.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.41%  '

$ws.Range("D39").Value = "'36.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("E40").Value = '  -2.72%  '

$ws.Range("D41").Value = "'0.778"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("D42").Value = "'273.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.18%  '

$ws.Range("E43").Value = '  -4.32%  '

$ws.Range("D44").Value = "'4.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.02%  '

$ws.Range("E45").Value = '  -1.65%  '

$ws.Range("D46").Value = "'0.0909"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.71%  '

$ws.Range("D47").Value = "'120.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.12%  '

$ws.Range("D48").Value = "'0.0485"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.20%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'17.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.51%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = "'0.0211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.69%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = "'16.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.65%  '
